# Fix heat rate modeling syntax
# Updates recalculated dispatch / state-of-charge / cost values across
# several sheets in the "Year 2" output workbook following a fix to the
# heat-rate modeling formula syntax.

$wb = $excel.ActiveWorkbook

# --- DG Dispatch: the single diesel-generator dispatch event shifts from
#     hour-column S (hour 17) to hour-column T (hour 18) on row 2.
$ws = $wb.Worksheets.Item("DG Dispatch")
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 2.884000000000018

# --- Costs and Revenues: totals recompute off the corrected dispatch.
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 76271.06239999998
$ws.Range("D2").Value = 9300.638068405266
$ws.Range("F2").Value = 44314.9292521668

# --- PV Dispatch
$ws = $wb.Worksheets.Item("PV Dispatch")
$ws.Range("H2").Value = 41.6
$ws.Range("I2").Value = 50.43636363636368
$ws.Range("L3").Value = 18.8531170288747
$ws.Range("M3").Value = 23.4
$ws.Range("O3").Value = 72.8
$ws.Range("R3").Value = 31.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 19.18312417100297
$ws.Range("P4").Value = 0

# --- Battery Input
$ws = $wb.Worksheets.Item("Battery Input")
$ws.Range("H2").Value = 28.6
$ws.Range("I2").Value = 19.23636363636368
$ws.Range("L3").Value = 18.8531170288747
$ws.Range("M3").Value = 0
$ws.Range("O3").Value = 72.8
$ws.Range("R3").Value = 31.2
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 19.18312417100297
$ws.Range("P4").Value = 0

# --- Battery Output
$ws = $wb.Worksheets.Item("Battery Output")
$ws.Range("S2").Value = 10.4
$ws.Range("T2").Value = 28.31599999999998

# --- State of Charge
$ws = $wb.Worksheets.Item("State of Charge")
$ws.Range("H2").Value = 170.784
$ws.Range("S2").Value = 637.4949494949495
$ws.Range("L3").Value = 333.5925858585859
$ws.Range("M3").Value = 333.5925858585859
$ws.Range("N3").Value = 333.5925858585859
$ws.Range("O3").Value = 405.664585858586
$ws.Range("P3").Value = 428.8305858585859
$ws.Range("Q3").Value = 454.570585858586
$ws.Range("J4").Value = 129.6
$ws.Range("K4").Value = 148.5912929292929
$ws.Range("L4").Value = 220.6632929292929
$ws.Range("M4").Value = 279.8652929292929
$ws.Range("N4").Value = 362.2332929292929
$ws.Range("O4").Value = 362.2332929292929
